$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Update the title text
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "4. Business Model Canvas – StudyMate AI"

# Update the content placeholder with the full Business Model Canvas breakdown
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$tr.Text = "1) Đối tác chính (Key Partners):`rOpenAI, HuggingFace, Google AI (AI/ML)`rTrường đại học, trung tâm giáo dục`rĐối tác thanh toán: Momo, ZaloPay, VNPay`rStartup EdTech, nhà xuất bản tài liệu`r2) Hoạt động chính (Key Activities):`rPhát triển & duy trì app (mobile/web)`rXây dựng/huấn luyện mô hình AI (tóm tắt, quiz, chatbot)`rMarketing online/offline tại trường học`rCSKH & hỗ trợ kỹ thuật`r3) Giá trị cốt lõi (Value Proposition):`rHọc thông minh hơn, tiết kiệm thời gian`rCá nhân hóa lộ trình, ôn tập hiệu quả`rTrợ lý ảo AI: tóm tắt, quiz, flashcard`rKhác biệt: tự động hóa – cá nhân hóa – tương tác như gia sư`r4) Quan hệ khách hàng (Customer Relationships):`rMiễn phí + nâng cấp Premium`rHỗ trợ chatbot 24/7, cộng đồng Facebook/Zalo`rGamification: tích điểm đổi thưởng`rEmail/SMS nhắc lịch học, deadline`r5) Phân khúc khách hàng (Customer Segments):`rSinh viên đại học, cao đẳng`rHọc sinh THPT chuẩn bị thi`rNgười đi làm muốn học thêm`r6) Kênh phân phối (Channels):`rApp Store, Google Play`rWebsite chính thức`rMXH: Facebook, TikTok, YouTube`rHợp tác CLB sinh viên, trung tâm gia sư`r7) Nguồn lực chính (Key Resources):`rĐội ngũ dev & chuyên gia AI`rHạ tầng cloud: AWS, GCP`rDữ liệu học tập (giáo trình, đề thi)`rVốn khởi nghiệp/đầu tư`r8) Cơ cấu chi phí (Cost Structure):`rPhát triển ứng dụng & server cloud`rNhân sự: dev, AI, marketing`rMarketing & quảng cáo`rBản quyền AI/API`r9) Dòng doanh thu (Revenue Streams):`rGói Premium: 99k/tháng (AI nâng cao, flashcard không giới hạn)`rQuảng cáo (phiên bản free)`rB2B: Giải pháp AI cho trường học/trung tâm`rKhóa học mini tích hợp trong app"

# Apply indent levels per paragraph (PowerPoint IndentLevel is 1-based: 1 = top level, 2 = sub-level)
# Only set IndentLevel on the sub-bullets; section headings keep the default top level.
$tr.Paragraphs(2).IndentLevel = 2
$tr.Paragraphs(3).IndentLevel = 2
$tr.Paragraphs(4).IndentLevel = 2
$tr.Paragraphs(5).IndentLevel = 2
$tr.Paragraphs(7).IndentLevel = 2
$tr.Paragraphs(8).IndentLevel = 2
$tr.Paragraphs(9).IndentLevel = 2
$tr.Paragraphs(10).IndentLevel = 2
$tr.Paragraphs(12).IndentLevel = 2
$tr.Paragraphs(13).IndentLevel = 2
$tr.Paragraphs(14).IndentLevel = 2
$tr.Paragraphs(15).IndentLevel = 2
$tr.Paragraphs(17).IndentLevel = 2
$tr.Paragraphs(18).IndentLevel = 2
$tr.Paragraphs(19).IndentLevel = 2
$tr.Paragraphs(20).IndentLevel = 2
$tr.Paragraphs(22).IndentLevel = 2
$tr.Paragraphs(23).IndentLevel = 2
$tr.Paragraphs(24).IndentLevel = 2
$tr.Paragraphs(26).IndentLevel = 2
$tr.Paragraphs(27).IndentLevel = 2
$tr.Paragraphs(28).IndentLevel = 2
$tr.Paragraphs(29).IndentLevel = 2
$tr.Paragraphs(31).IndentLevel = 2
$tr.Paragraphs(32).IndentLevel = 2
$tr.Paragraphs(33).IndentLevel = 2
$tr.Paragraphs(34).IndentLevel = 2
$tr.Paragraphs(36).IndentLevel = 2
$tr.Paragraphs(37).IndentLevel = 2
$tr.Paragraphs(38).IndentLevel = 2
$tr.Paragraphs(39).IndentLevel = 2
$tr.Paragraphs(41).IndentLevel = 2
$tr.Paragraphs(42).IndentLevel = 2
$tr.Paragraphs(43).IndentLevel = 2
$tr.Paragraphs(44).IndentLevel = 2
